$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3872.16
$ws.Range("I28").Value = 1190.7333
$ws.Range("K28").Value = 1190.7333
$ws.Range("M28").Value = -705.7333000000001
$ws.Range("H52").Value = 1293
$ws.Range("I52").Value = 1100
$ws.Range("J52").Value = 1389.5
$ws.Range("K52").Value = 3300
$ws.Range("L52").Value = 4168.5
$ws.Range("M52").Value = -3140
$ws.Range("N52").Value = -4488.5
$ws.Range("H62").Value = 2889.889
$ws.Range("I62").Value = 2318.1428
$ws.Range("K62").Value = 2318.1428
$ws.Range("M62").Value = -1694.1428
$ws.Range("H65").Value = 2889.889
$ws.Range("I65").Value = 2318.1428
$ws.Range("K65").Value = 11590.714
$ws.Range("M65").Value = -8470.714
$ws.Range("H98").Value = 1435.4062
$ws.Range("I98").Value = 1244.8462
$ws.Range("K98").Value = 1244.8462
$ws.Range("M98").Value = 253.1538
$ws.Range("H99").Value = 3067.5833
$ws.Range("I99").Value = 164.6
$ws.Range("J99").Value = 5141.143
$ws.Range("K99").Value = 493.8
$ws.Range("L99").Value = 15423.429
$ws.Range("M99").Value = 1004.2
$ws.Range("N99").Value = -18419.429
$ws.Range("H113").Value = 8899.5
$ws.Range("I113").Value = 8899.5
$ws.Range("K113").Value = 8899.5
$ws.Range("M113").Value = -5645.5
$ws.Range("H122").Value = 1435.4062
$ws.Range("I122").Value = 1244.8462
$ws.Range("K122").Value = 3734.5386
$ws.Range("M122").Value = -1284.5386
$ws.Range("H127").Value = 2911.2856
$ws.Range("I127").Value = 1349.75
$ws.Range("J127").Value = 4993.3335
$ws.Range("K127").Value = 4049.25
$ws.Range("L127").Value = 14980.0005
$ws.Range("M127").Value = 910.75
$ws.Range("N127").Value = -24900.0005
$ws.Range("H132").Value = 4325.1313
$ws.Range("I132").Value = 2257.276
$ws.Range("J132").Value = 10988.223
$ws.Range("K132").Value = 6771.828
$ws.Range("L132").Value = 32964.669
$ws.Range("M132").Value = -4241.828
$ws.Range("N132").Value = -38024.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""
$ws.Range("H32").Value = 1834.56
$ws.Range("I32").Value = 1624.5217
$ws.Range("K32").Value = 1624.5217
$ws.Range("M32").Value = -1337.5217
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1594
$ws.Range("H74").Value = 2277529.5
$ws.Range("I74").Value = 3575118.5
$ws.Range("K74").Value = 3575118.5
$ws.Range("M74").Value = -3574244.5
$ws.Range("H77").Value = 2277529.5
$ws.Range("I77").Value = 3575118.5
$ws.Range("K77").Value = 17875592.5
$ws.Range("M77").Value = -17871224.5
$ws.Range("H102").Value = 2050
$ws.Range("I102").Value = 1863.9286
$ws.Range("K102").Value = 1863.9286
$ws.Range("M102").Value = -241.9286
$ws.Range("H132").Value = 3337987
$ws.Range("I132").Value = 4814.241
$ws.Range("K132").Value = 14442.723
$ws.Range("M132").Value = -11912.723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3957.6
$ws.Range("I20").Value = 3697
$ws.Range("K20").Value = 3697
$ws.Range("M20").Value = -3450
$ws.Range("H80").Value = 1182.7778
$ws.Range("J80").Value = 1374.5
$ws.Range("L80").Value = 1374.5
$ws.Range("N80").Value = -3370.5
$ws.Range("H83").Value = 1182.7778
$ws.Range("J83").Value = 1374.5
$ws.Range("L83").Value = 6872.5
$ws.Range("N83").Value = -16856.5
$ws.Range("H86").Value = 21655.258
$ws.Range("I86").Value = 38808.223
$ws.Range("K86").Value = 38808.223
$ws.Range("M86").Value = -37685.223
$ws.Range("H89").Value = 21655.258
$ws.Range("I89").Value = 38808.223
$ws.Range("K89").Value = 194041.115
$ws.Range("M89").Value = -188425.115
$ws.Range("H94").Value = 1326.3954
$ws.Range("I94").Value = 982.63336
$ws.Range("J94").Value = 2119.6924
$ws.Range("K94").Value = 982.63336
$ws.Range("L94").Value = 2119.6924
$ws.Range("M94").Value = -531.63336
$ws.Range("N94").Value = -3021.6924
$ws.Range("H99").Value = 3038
$ws.Range("I99").Value = 3250.7778
$ws.Range("J99").Value = 2399.6667
$ws.Range("K99").Value = 3250.7778
$ws.Range("L99").Value = 2399.6667
$ws.Range("M99").Value = -1752.7778
$ws.Range("N99").Value = -5395.6667
$ws.Range("H107").Value = 8537.625
$ws.Range("I107").Value = 9414.666999999999
$ws.Range("K107").Value = 9414.666999999999
$ws.Range("M107").Value = -7494.666999999999
$ws.Range("H113").Value = 15020
$ws.Range("I113").Value = 15020
$ws.Range("K113").Value = 15020
$ws.Range("M113").Value = -12850
$ws.Range("H134").Value = 7144923
$ws.Range("I134").Value = 2077.0833
$ws.Range("K134").Value = 6231.249899999999
$ws.Range("M134").Value = -3696.249899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3264.375
$ws.Range("I58").Value = 3264.375
$ws.Range("K58").Value = 3264.375
$ws.Range("M58").Value = -3061.375
$ws.Range("H99").Value = 40612
$ws.Range("I99").Value = 12998.75
$ws.Range("K99").Value = 12998.75
$ws.Range("M99").Value = -11500.75
$ws.Range("H107").Value = 1492.7142
$ws.Range("I107").Value = 1299.1578
$ws.Range("K107").Value = 1299.1578
$ws.Range("M107").Value = 620.8422
$ws.Range("H126").Value = 40612
$ws.Range("I126").Value = 12998.75
$ws.Range("K126").Value = 38996.25
$ws.Range("M126").Value = -36526.25
$ws.Range("H134").Value = 2999.4119
$ws.Range("I134").Value = 2999.4119
$ws.Range("K134").Value = 8998.235700000001
$ws.Range("M134").Value = -6463.235700000001
$ws.Range("H136").Value = 3264.375
$ws.Range("I136").Value = 3264.375
$ws.Range("K136").Value = 9793.125
$ws.Range("M136").Value = -7243.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 29444.334
$ws.Range("I114").Value = 30000
$ws.Range("J114").Value = 29166.5
$ws.Range("K114").Value = 90000
$ws.Range("L114").Value = 87499.5
$ws.Range("M114").Value = -86746
$ws.Range("N114").Value = -94007.5
$ws.Range("H129").Value = 8015.778
$ws.Range("I129").Value = 4974.2
$ws.Range("K129").Value = 14922.6
$ws.Range("M129").Value = -9922.599999999999
$ws.Range("H139").Value = 6671.8335
$ws.Range("I139").Value = 1339.6
$ws.Range("K139").Value = 4018.8
$ws.Range("M139").Value = 1121.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 635.2963
$ws.Range("I97").Value = 508.4375
$ws.Range("J97").Value = 819.8182
$ws.Range("K97").Value = 508.4375
$ws.Range("L97").Value = 819.8182
$ws.Range("M97").Value = -12.4375
$ws.Range("N97").Value = -1811.8182
$ws.Range("H113").Value = 807215.25
$ws.Range("I113").Value = 2102.5625
$ws.Range("K113").Value = 2102.5625
$ws.Range("M113").Value = 67.4375
$ws.Range("H126").Value = 2486.1
$ws.Range("I126").Value = 2724.6667
$ws.Range("J126").Value = 2128.25
$ws.Range("K126").Value = 8174.000100000001
$ws.Range("L126").Value = 6384.75
$ws.Range("M126").Value = -5704.000100000001
$ws.Range("N126").Value = -11324.75
$ws.Range("I132").Value = 5498.5
$ws.Range("K132").Value = 16495.5
$ws.Range("M132").Value = -13965.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2838.0908
$ws.Range("I46").Value = 2470
$ws.Range("J46").Value = 3144.8333
$ws.Range("K46").Value = 2470
$ws.Range("L46").Value = 3144.8333
$ws.Range("M46").Value = -2282
$ws.Range("N46").Value = -3520.8333
$ws.Range("H82").Value = 4870.8667
$ws.Range("I82").Value = 4300.4
$ws.Range("J82").Value = 5156.1
$ws.Range("K82").Value = 4300.4
$ws.Range("L82").Value = 5156.1
$ws.Range("M82").Value = -3939.4
$ws.Range("N82").Value = -5878.1
$ws.Range("H85").Value = 4870.8667
$ws.Range("I85").Value = 4300.4
$ws.Range("J85").Value = 5156.1
$ws.Range("K85").Value = 4300.4
$ws.Range("L85").Value = 5156.1
$ws.Range("M85").Value = -3052.4
$ws.Range("N85").Value = -7652.1
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H122").Value = 3767.5
$ws.Range("I122").Value = 3377.3547
$ws.Range("J122").Value = 7799
$ws.Range("K122").Value = 10132.0641
$ws.Range("L122").Value = 23397
$ws.Range("M122").Value = -7682.0641
$ws.Range("N122").Value = -28297
$ws.Range("H136").Value = 3650.6
$ws.Range("I136").Value = 2945.111
$ws.Range("K136").Value = 8835.332999999999
$ws.Range("M136").Value = -6285.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H126").Value = 6432.067
$ws.Range("I126").Value = 6666.75
$ws.Range("K126").Value = 20000.25
$ws.Range("M126").Value = -17530.25
